# Generate Report for Handoff
# Refresh the "Latest Handoff Date(time)" stamps for the files that are
# part of the new handoff batch (status "Handback transform failed" and
# "Ready for handoff") on all three report sheets. Rows that are already
# "Handed back: in sync with en-US" or still "In Translation" are left
# untouched.

$wb = $excel.ActiveWorkbook

# Rows (1-based, matching worksheet row numbers) whose handoff
# date/datetime gets refreshed - 4, 6, 7, 8, 9, 10.
$rows = @(4, 6, 7, 8, 9, 10)

# --- Overview sheet: column D = "Latest Handoff Date" ---
$wsOverview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $wsOverview.Cells.Item($r, 4).Value = "2016-03-23 07:36:08"
}

# --- zh-cn sheet: column E = "Latest Handoff Datetime" ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $wsZhCn.Cells.Item($r, 5).Value = "2016-03-23 07:35:59"
}

# --- de-de sheet: column E = "Latest Handoff Datetime" ---
$wsDeDe = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $wsDeDe.Cells.Item($r, 5).Value = "2016-03-23 07:36:08"
}
